$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to text format first so numeric-looking
# strings like "315.80" or "1.002" are stored as text, matching the
# original inlineStr cell contents instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '24.915.87'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.703.78'
$ws.Range("E4").Value = '  -0.53%  '
$ws.Range("D5").Value = '315.80'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("D7").Value = '0.4060'
$ws.Range("E7").Value = '  +2.80%  '
$ws.Range("D8").Value = '0.4063'
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").Value = '1.002'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("D10").Value = '53.71'
$ws.Range("E10").Value = '  +2.40%  '
$ws.Range("D11").Value = '1.468'
$ws.Range("E11").Value = '  -3.14%  '
$ws.Range("D12").Value = '0.08813'
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").Value = '25.81'
$ws.Range("E13").Value = '  +4.14%  '
$ws.Range("D14").Value = '7.503'
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").Value = '8.051'
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("D17").Value = '1.738.27'
$ws.Range("E17").Value = '  +2.24%  '
$ws.Range("D18").Value = '96.57'
$ws.Range("D19").Value = '0.07168'
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("D20").Value = '20.98'
$ws.Range("E20").Value = '  +4.84%  '
$ws.Range("D21").Value = '7.237'
$ws.Range("E21").Value = '  -2.18%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.94%  '
$ws.Range("D23").Value = '14.60'
$ws.Range("E23").Value = '  +1.48%  '
$ws.Range("D24").Value = '24.941.95'
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").Value = '2.325'
$ws.Range("E25").Value = '  -1.13%  '
$ws.Range("D26").Value = '6.840'
$ws.Range("E26").Value = '  +30.85%  '
$ws.Range("D27").Value = '2.885'
$ws.Range("E27").Value = '  -5.62%  '
$ws.Range("D28").Value = '23.05'
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("D29").Value = '164.80'
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("D30").Value = '145.25'
$ws.Range("E30").Value = '  +4.82%  '
$ws.Range("D31").Value = '8.221'
$ws.Range("E31").Value = '  -5.72%  '
$ws.Range("E32").Value = '  +14.16%  '
$ws.Range("D33").Value = '1.913.45'
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("D34").Value = '0.08799'
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("D35").Value = '0.03199'
$ws.Range("E35").Value = '  +10.26%  '
$ws.Range("D36").Value = '7.319'
$ws.Range("E36").Value = '  -5.14%  '
$ws.Range("D37").Value = '1.016'
$ws.Range("E37").Value = '  -3.38%  '
$ws.Range("E38").Value = '  +3.54%  '
$ws.Range("D39").Value = '0.8469'
$ws.Range("E39").Value = '  +7.47%  '
$ws.Range("D40").Value = '10.91'
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("D41").Value = '0.09419'
$ws.Range("E41").Value = '  +2.84%  '
$ws.Range("D42").Value = '14.05'
$ws.Range("E42").Value = '  -2.71%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.471'
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '17.76'
$ws.Range("E44").Value = '  +6.13%  '
$ws.Range("D45").Value = '2.727'
$ws.Range("E45").Value = '  +4.26%  '
$ws.Range("D46").Value = '0.7436'
$ws.Range("E46").Value = '  +2.41%  '
$ws.Range("D47").Value = '4.234'
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("D48").Value = '1.395'
$ws.Range("E48").Value = '  +4.55%  '
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("D50").Value = '142.20'
$ws.Range("E50").Value = '  +1.25%  '
$ws.Range("D51").Value = '0.08358'
$ws.Range("E51").Value = '  +3.92%  '

# Restore default (unstyled) formatting on the touched columns so we do not
# leave a stray text-format style on cells that did not have one originally.
$ws.Range("D2:E51").Style = "Normal"
